$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear columns B and C entirely (they are no longer used)
$ws.Range("B:C").Clear()

# Clear rows 5 through 13 (old Manhattan/Minkowski/Supremum rows no longer needed)
$ws.Range("A5:A13").Clear()

# Row 1: header row collapses to a single numeric placeholder cell, keep its style
$ws.Range("A1").Value = 0

# Rows 2-4: new car names (Supremum-only result set)
$ws.Range("A2").Value = "Suzuki Ertiga"
$ws.Range("A3").Value = "Xpander"
$ws.Range("A4").Value = "Livina"
